# Introduction to Mockito presentation updates
#
# Slide 15 ("Using the InlineMockMaker") goes from two empty
# placeholders (Title + Content) to a filled-in title and four
# bullet paragraphs describing the InlineMockMaker / MockMaker API,
# with the mockito-core / mockito-inline artifacts at indent level 2.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)

# ---------------------------------------------------------------
# Shape 1 (Title): "Using the InlineMockMaker"
# ---------------------------------------------------------------
$titleShape = $s.Shapes.Item(1)
$titleTr = $titleShape.TextFrame.TextRange
$titleTr.Text = "Using the InlineMockMaker"
$titleShape.TextFrame.TextRange.Characters(11, 15).Font.Name = "Courier"

# ---------------------------------------------------------------
# Shape 2 (Content placeholder): 4 paragraphs
# ---------------------------------------------------------------
$bodyShape = $s.Shapes.Item(2)
$bodyTr = $bodyShape.TextFrame.TextRange

# Paragraph 1
$bodyTr.Text = "The InlineMockMaker is the MockMaker API implementation that provides static method and constructor mocking."
$bodyShape.TextFrame.TextRange.Characters(5, 15).Font.Name = "Courier"

# Paragraph 2
$bodyShape.TextFrame.TextRange.InsertAfter("`rMockMaker is experimental.")

# Paragraph 3 (indent level 2 -> lvl="1")
$bodyShape.TextFrame.TextRange.InsertAfter("`rFunctionality is in the mockito-core dependency artifact, but is turned off by default. Requires configuring the MockMaker extension file.")

# Paragraph 4 (indent level 2 -> lvl="1")
$bodyShape.TextFrame.TextRange.InsertAfter("`rThe mockito-inline dependency artifact preconfigures mock maker for inline mock making. Thus no need for the MockMaker extension file.")

# Indent levels for paragraphs 3 & 4
$bodyShape.TextFrame.TextRange.Paragraphs(3).IndentLevel = 2
$bodyShape.TextFrame.TextRange.Paragraphs(4).IndentLevel = 2

# Bold "mockito" / "-core" in paragraph 3
$bodyShape.TextFrame.TextRange.Paragraphs(3).Characters(25, 7).Font.Bold = $true
$bodyShape.TextFrame.TextRange.Paragraphs(3).Characters(32, 5).Font.Bold = $true

# Bold "mockito" / "-inline" in paragraph 4
$bodyShape.TextFrame.TextRange.Paragraphs(4).Characters(5, 7).Font.Bold = $true
$bodyShape.TextFrame.TextRange.Paragraphs(4).Characters(12, 7).Font.Bold = $true
